$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I35").Value = 'sv'
$ws.Range("J35").Value = 'Statement-opinion'
$ws.Range("I36").Value = 'ba'
$ws.Range("J36").Value = 'Appreciation'
$ws.Range("I39").Value = 'b'
$ws.Range("J39").Value = 'Acknowledge (Backchannel)'
$ws.Range("I47").Value = 'b'
$ws.Range("J47").Value = 'Acknowledge (Backchannel)'
$ws.Range("I48").Value = 'sd'
$ws.Range("J48").Value = 'Statement-non-opinion'
$ws.Range("I51").Value = 'sd'
$ws.Range("J51").Value = 'Statement-non-opinion'
$ws.Range("I52").Value = 'sd'
$ws.Range("J52").Value = 'Statement-non-opinion'
$ws.Range("I63").Value = 'sd'
$ws.Range("J63").Value = 'Statement-non-opinion'
$ws.Range("I84").Value = 'b'
$ws.Range("J84").Value = 'Acknowledge (Backchannel)'
$ws.Range("I95").Value = 'sd'
$ws.Range("J95").Value = 'Statement-non-opinion'
$ws.Range("I111").Value = 'sd'
$ws.Range("J111").Value = 'Statement-non-opinion'
$ws.Range("I170").Value = 'aa'
$ws.Range("J170").Value = 'Agree/Accept'
$ws.Range("I171").Value = 'aa'
$ws.Range("J171").Value = 'Agree/Accept'
$ws.Range("I172").Value = 'sd'
$ws.Range("J172").Value = 'Statement-non-opinion'
$ws.Range("I187").Value = 'ba'
$ws.Range("J187").Value = 'Appreciation'
$ws.Range("I236").Value = 'aa'
$ws.Range("J236").Value = 'Agree/Accept'
$ws.Range("I251").Value = 'sd'
$ws.Range("J251").Value = 'Statement-non-opinion'
$ws.Range("I264").Value = 'aa'
$ws.Range("J264").Value = 'Agree/Accept'
$ws.Range("I290").Value = 'ba'
$ws.Range("J290").Value = 'Appreciation'
$ws.Range("I292").Value = 'aa'
$ws.Range("J292").Value = 'Agree/Accept'
$ws.Range("I298").Value = 'sv'
$ws.Range("J298").Value = 'Statement-opinion'
$ws.Range("I299").Value = 'sd'
$ws.Range("J299").Value = 'Statement-non-opinion'
$ws.Range("I303").Value = 'sd'
$ws.Range("J303").Value = 'Statement-non-opinion'
$ws.Range("I307").Value = 'sd'
$ws.Range("J307").Value = 'Statement-non-opinion'
$ws.Range("I312").Value = 'sv'
$ws.Range("J312").Value = 'Statement-opinion'
$ws.Range("I325").Value = 'aa'
$ws.Range("J325").Value = 'Agree/Accept'
$ws.Range("I331").Value = 'sd'
$ws.Range("J331").Value = 'Statement-non-opinion'
$ws.Range("I350").Value = 'aa'
$ws.Range("J350").Value = 'Agree/Accept'
$ws.Range("I353").Value = 'ba'
$ws.Range("J353").Value = 'Appreciation'
$ws.Range("I393").Value = 'sv'
$ws.Range("J393").Value = 'Statement-opinion'
$ws.Range("I406").Value = 'sd'
$ws.Range("J406").Value = 'Statement-non-opinion'
